# "added the total score function"
#
# The sentiment-analysis output sheet gains several new metric columns.
# The existing "Method" column (last-used column, narrow width) is moved
# to become the last column (M), and four new percentage columns
# (posWordPercentage, negWordPercentage, posPhrasePercentage,
# negPhrasePercentage) are inserted right after "Date", pushing the other
# metric columns (ElapsedMs, wordCount, sentenceCount, posWordCount,
# negWordCount, positivePhraseCount, negativePhraseCount) further right.
#
# New header order (row 1):
#   A Date
#   B posWordPercentage
#   C negWordPercentage
#   D posPhrasePercentage
#   E negPhrasePercentage
#   F ElapsedMs
#   G wordCount
#   H sentenceCount
#   I posWordCount
#   J negWordCount
#   K positivePhraseCount
#   L negativePhraseCount
#   M Method
#
# Row 2 is updated with a new sample result (and a later "Date" timestamp)
# computed under the new, richer scoring function.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header row (row 1) ----
$ws.Range("B1").Value = "posWordPercentage"
$ws.Range("C1").Value = "negWordPercentage"
$ws.Range("D1").Value = "posPhrasePercentage"
$ws.Range("E1").Value = "negPhrasePercentage"
$ws.Range("F1").Value = "ElapsedMs"
$ws.Range("G1").Value = "wordCount"
$ws.Range("H1").Value = "sentenceCount"
$ws.Range("I1").Value = "posWordCount"
$ws.Range("J1").Value = "negWordCount"
$ws.Range("K1").Value = "positivePhraseCount"
$ws.Range("L1").Value = "negativePhraseCount"
$ws.Range("M1").Value = "Method"

# ---- Data row (row 2) ----
$ws.Range("A2").Value = 42605.455104166664
$ws.Range("B2").Value = 73
$ws.Range("C2").Value = 26
$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 100
$ws.Range("F2").Value = 7474
$ws.Range("G2").Value = 3826
$ws.Range("H2").Value = 202
$ws.Range("I2").Value = 73
$ws.Range("J2").Value = 26
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 5
$ws.Range("M2").Value = "Named"
